# daily auto push: 2026-01-22 02:33 UTC
#
# A new sample row (2026/01/22, 木, 10, 166) was recorded between the
# existing "2026/01/22" rows (row 701) and the "2026/12/29" rows
# (old row 702). Insert one row at row 702 to push everything else down
# by one and fill in the new reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 702..743 down to 703..744, creating a blank row 702.
$ws.Rows.Item(702).Insert()

# The date column stores dates as plain text (e.g. "2026/01/22"), not
# real date values, so force a text number format before assigning the
# value to stop Excel from auto-converting it to a date serial number.
$ws.Range("A702").NumberFormat = "@"
$ws.Range("A702").Value = "2026/01/22"
# Restore the default ("Normal") cell style so the new cell matches the
# unstyled data cells around it instead of keeping the temporary text
# format applied above.
$ws.Range("A702").Style = "Normal"

$ws.Range("B702").Value = "木"
$ws.Range("C702").Value = 10
$ws.Range("D702").Value = 166
